$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4365890793924621
$ws.Range("B3").Value = 0.41495687484462
$ws.Range("B4").Value = 0.402163730678609
$ws.Range("B5").Value = 0.3970726832776847
$ws.Range("B6").Value = 0.3962346820407845
$ws.Range("B7").Value = 0.4020945769558466
$ws.Range("B8").Value = 0.4290283276575622
$ws.Range("B9").Value = 0.4857633054946575
$ws.Range("B10").Value = 0.5298922905263908
$ws.Range("B11").Value = 0.5505117540847664
$ws.Range("B12").Value = 0.5583991422754195
$ws.Range("B13").Value = 0.5566969134648332
$ws.Range("B14").Value = 0.5511590609695816
$ws.Range("B15").Value = 0.5477773128518209
$ws.Range("B16").Value = 0.5285558009489932
$ws.Range("B17").Value = 0.516904224727142
$ws.Range("B18").Value = 0.5102537874046504
$ws.Range("B19").Value = 0.5080108378277544
$ws.Range("B20").Value = 0.5181392453346234
$ws.Range("B21").Value = 0.552783504702802
$ws.Range("B22").Value = 0.5758878521774307
$ws.Range("B23").Value = 0.5635140291904861
$ws.Range("B24").Value = 0.5175807424232914
$ws.Range("B25").Value = 0.4699896866073914

$ws.Range("C2").Value = 0.0229848125899963
$ws.Range("C3").Value = 0.02054100797682423
$ws.Range("C4").Value = 0.01908389746531469
$ws.Range("C5").Value = 0.01850080003072918
$ws.Range("C6").Value = 0.01840461606765587
$ws.Range("C7").Value = 0.01907599063510901
$ws.Range("C8").Value = 0.02213304509371028
$ws.Range("C9").Value = 0.02848278974953189
$ws.Range("C10").Value = 0.03337945442947898
$ws.Range("C11").Value = 0.03566069354960177
$ws.Range("C12").Value = 0.03653253243929555
$ws.Range("C13").Value = 0.03634440713275922
$ws.Range("C14").Value = 0.035732258773578
$ws.Range("C15").Value = 0.03535834759371426
$ws.Range("C16").Value = 0.03323147307578722
$ws.Range("C17").Value = 0.03194064625461124
$ws.Range("C18").Value = 0.03120323769613265
$ws.Range("C19").Value = 0.030954421807877
$ws.Range("C20").Value = 0.03207753340335273
$ws.Range("C21").Value = 0.03591184284987037
$ws.Range("C22").Value = 0.03846444013483108
$ws.Range("C23").Value = 0.03709771577422316
$ws.Range("C24").Value = 0.03201563208817504
$ws.Range("C25").Value = 0.0267254691854788

$ws.Range("D2").Value = 0.2179683731167756
$ws.Range("D3").Value = 0.2057779472138606
$ws.Range("D4").Value = 0.1982317192729255
$ws.Range("D5").Value = 0.195140718049629
$ws.Range("D6").Value = 0.1946264891473959
$ws.Range("D7").Value = 0.1981900977190776
$ws.Range("D8").Value = 0.2137775550301342
$ws.Range("D9").Value = 0.2438798400194315
$ws.Range("D10").Value = 0.265743963899439
$ws.Range("D11").Value = 0.2756426529658995
$ws.Range("D12").Value = 0.2793847332979169
$ws.Range("D13").Value = 0.2785790835772843
$ws.Range("D14").Value = 0.2759506404388219
$ws.Range("D15").Value = 0.2743398319631751
$ws.Range("D16").Value = 0.265096148147336
$ws.Range("D17").Value = 0.2594136443586592
$ws.Range("D18").Value = 0.2561407270967351
$ws.Range("D19").Value = 0.2550317876130919
$ws.Range("D20").Value = 0.2600190179680624
$ws.Range("D21").Value = 0.2767228460881199
$ws.Range("D22").Value = 0.2876030434817949
$ws.Range("D23").Value = 0.2817992779401379
$ws.Range("D24").Value = 0.2597453471911422
$ws.Range("D25").Value = 0.235782901811902

$ws.Range("E2").Value = 0.1375878983431775
$ws.Range("E3").Value = 0.1255670726398108
$ws.Range("E4").Value = 0.1182492430580666
$ws.Range("E5").Value = 0.1152826547156707
$ws.Range("E6").Value = 0.1147909799882854
$ws.Range("E7").Value = 0.1182091723704985
$ws.Range("E8").Value = 0.1334298213663629
$ws.Range("E9").Value = 0.1637951062884611
$ws.Range("E10").Value = 0.1864469207310222
$ws.Range("E11").Value = 0.1968321778962974
$ws.Range("E12").Value = 0.200776875403136
$ws.Range("E13").Value = 0.1999267728851351
$ws.Range("E14").Value = 0.1971564671884423
$ws.Range("E15").Value = 0.1954611531454944
$ws.Range("E16").Value = 0.1857698805716268
$ws.Range("E17").Value = 0.1798455863790238
$ws.Range("E18").Value = 0.1764456730249222
$ws.Range("E19").Value = 0.1752958111804617
$ws.Range("E20").Value = 0.1804754498020174
$ws.Range("E21").Value = 0.1979698435059376
$ws.Range("E22").Value = 0.2094737290249498
$ws.Range("E23").Value = 0.2033273229360191
$ws.Range("E24").Value = 0.1801906697584599
$ws.Range("E25").Value = 0.1555223178531051

$ws.Range("F2").Value = 3.77879553085458
$ws.Range("F3").Value = 3.557895458523234
$ws.Range("F4").Value = 3.422424573983392
$ws.Range("F5").Value = 3.367255022453122
$ws.Range("F6").Value = 3.358096186686765
$ws.Range("F7").Value = 3.421680399169105
$ws.Range("F8").Value = 3.702592361287145
$ws.Range("F9").Value = 4.255001867998317
$ws.Range("F10").Value = 4.662184838814426
$ws.Range("F11").Value = 4.847798002161312
$ws.Range("F12").Value = 4.918146339120653
$ws.Range("F13").Value = 4.902992779721785
$ws.Range("F14").Value = 4.853584354661962
$ws.Range("F15").Value = 4.823328332074823
$ws.Range("F16").Value = 4.650062769241345
$ws.Range("F17").Value = 4.543872093339985
$ws.Range("F18").Value = 4.482829796502671
$ws.Range("F19").Value = 4.46216790109878
$ws.Range("F20").Value = 4.55517253028242
$ws.Range("F21").Value = 4.868095116235907
$ws.Range("F22").Value = 5.072964301113245
$ws.Range("F23").Value = 4.963587331918802
$ws.Range("F24").Value = 4.55006357780934
$ws.Range("F25").Value = 4.105355574994093

$ws.Range("K2").Value = 0.3818434170172509
$ws.Range("K3").Value = 0.3602888624461684
$ws.Range("K4").Value = 0.3475414940870678
$ws.Range("K5").Value = 0.3424684131769311
$ws.Range("K6").Value = 0.3416333476347404
$ws.Range("K7").Value = 0.3474725857253276
$ws.Range("K8").Value = 0.3743097075395667
$ws.Range("K9").Value = 0.4308494698544223
$ws.Range("K10").Value = 0.4748443858160272
$ws.Range("K11").Value = 0.4954076795042681
$ws.Range("K12").Value = 0.5032747244793825
$ws.Range("K13").Value = 0.5015768319546794
$ws.Range("K14").Value = 0.4960532928090231
$ws.Range("K15").Value = 0.492680439761358
$ws.Range("K16").Value = 0.4735116839950706
$ws.Range("K17").Value = 0.4618938719459322
$ws.Range("K18").Value = 0.4552632942943546
$ws.Range("K19").Value = 0.4530271372381094
$ws.Range("K20").Value = 0.4631252511346986
$ws.Range("K21").Value = 0.4976735053242862
$ws.Range("K22").Value = 0.5207206222670777
$ws.Range("K23").Value = 0.5083767616616797
$ws.Range("K24").Value = 0.4625683931307378
$ws.Range("K25").Value = 0.4151279556318173

$ws.Range("M2").Value = 0.2546306249913073
$ws.Range("M3").Value = 0.2359885170281188
$ws.Range("M4").Value = 0.2247481408713696
$ws.Range("M5").Value = 0.2202187676389897
$ws.Range("M6").Value = 0.2194697414104425
$ws.Range("M7").Value = 0.2246868497770649
$ws.Range("M8").Value = 0.2481597004513887
$ws.Range("M9").Value = 0.2958547618691796
$ws.Range("M10").Value = 0.3319579783199273
$ws.Range("M11").Value = 0.3486232691084155
$ws.Range("M12").Value = 0.3549695319698927
$ws.Range("M13").Value = 0.3536011612638177
$ws.Range("M14").Value = 0.3491446647835801
$ws.Range("M15").Value = 0.346419573843022
$ws.Range("M16").Value = 0.3308737879318002
$ws.Range("M17").Value = 0.3213993759142468
$ws.Range("M18").Value = 0.3159726564784791
$ws.Range("M19").Value = 0.314139141723274
$ws.Range("M20").Value = 0.3224055875299356
$ws.Range("M21").Value = 0.3504526771804279
$ws.Range("M22").Value = 0.3689902107869543
$ws.Range("M23").Value = 0.3590771875059033
$ws.Range("M24").Value = 0.3219506162498007
$ws.Range("M25").Value = 0.2827689630748225

